$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (values stored as text, same as original inline strings)
$ws.Range("B2").Value = "'40"
$ws.Range("C2").Value = "'40"
$ws.Range("D2").Value = "'39"
$ws.Range("E2").Value = "'39"
$ws.Range("F2").Value = "'42"
$ws.Range("G2").Value = "'41"
$ws.Range("H2").Value = "'41"
$ws.Range("I2").Value = "'41"
$ws.Range("J2").Value = "'39"
$ws.Range("K2").Value = "'41"
$ws.Range("L2").Value = "'40"
$ws.Range("M2").Value = "'42"
$ws.Range("N2").Value = "'540"
$ws.Range("O2").Value = "'475"

# Row 3 (values stored as text, same as original inline strings)
$ws.Range("B3").Value = "'19"
$ws.Range("C3").Value = "'19"
$ws.Range("D3").Value = "'18"
$ws.Range("E3").Value = "'19"
$ws.Range("F3").Value = "'22"
$ws.Range("G3").Value = "'20"
$ws.Range("H3").Value = "'19"
$ws.Range("I3").Value = "'20"
$ws.Range("J3").Value = "'20"
$ws.Range("K3").Value = "'19"
$ws.Range("L3").Value = "'19"
$ws.Range("M3").Value = "'19"
$ws.Range("N3").Value = "'300"
$ws.Range("O3").Value = "'229"
